$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PAYC")

# Row 6: "Change in inventories"
$ws.Range("B6").Value = -41000.0
$ws.Range("C6").Value = -190000.0
$ws.Range("D6").Value = -588000.0
$ws.Range("E6").Value = -138000.0
$ws.Range("F6").Value = -403000.0

# Row 8: "Change in payables and accrued liability"
$ws.Range("B8").Value = 26000000.0
$ws.Range("C8").Value = 24000000.0
$ws.Range("D8").Value = 20000000.0
$ws.Range("E8").Value = 16298000.0
$ws.Range("F8").Value = 15868000.0
